$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.353.88"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.862.90"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4770"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2754"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "1.868.53"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6319"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.04%  "
$ws.Range("D16").Value = "30.298.01"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007378"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "2.099.07"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.093"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.3944"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.006"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.291"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.862"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.382"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.95%  "
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04914"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7240"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01920"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.79%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9046"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "105.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4105"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.565"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.063"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "61.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.806"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.404"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.37%  "
